$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lowcode.dropDownSymbols: rows A2:F25 were shuffled into the correct
# symbol/reel-weight order. Rewrite each row explicitly to match.
$rows = @(
    @(2,201,9,30,15,45,30),
    @(3,1001,18,30,75,60,72),
    @(4,301,6,45,30,60,45),
    @(5,401,9,48,67,75,45),
    @(6,601,9,60,67,60,42),
    @(7,901,16,15,45,60,60),
    @(8,501,9,52,30,75,45),
    @(9,1201,2,10,10,10,10),
    @(10,101,9,30,15,60,15),
    @(11,902,1,0,0,0,0),
    @(12,801,3,67,65,52,45),
    @(13,701,3,90,45,97,15),
    @(14,1202,2,10,10,10,10),
    @(15,1203,3,15,15,15,15),
    @(16,1101,0,15,30,30,0),
    @(17,1,0,2,2,2,2),
    @(18,802,0,4,5,4,0),
    @(19,2,0,2,2,2,2),
    @(20,3,0,3,3,3,3),
    @(21,502,0,4,0,0,0),
    @(22,402,0,0,4,0,0),
    @(23,602,0,0,4,0,9),
    @(24,702,0,0,0,4,0),
    @(25,1002,0,0,0,0,9)
)

foreach ($row in $rows) {
    $rowNum = $row[0]
    $ws.Cells.Item($rowNum, 1).Value = $row[1]
    $ws.Cells.Item($rowNum, 2).Value = $row[2]
    $ws.Cells.Item($rowNum, 3).Value = $row[3]
    $ws.Cells.Item($rowNum, 4).Value = $row[4]
    $ws.Cells.Item($rowNum, 5).Value = $row[5]
    $ws.Cells.Item($rowNum, 6).Value = $row[6]
}

Write-Host "Applied row permutation to A2:F25"
